$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on numeric-looking Price cells so they stay literal strings
# (matches the source workbook, where these are inline strings, not numbers).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.268.87"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "1.922.58"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").Value = "318.04"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "0.4866"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").Value = "0.3851"
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").Value = "0.07405"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "0.9435"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "20.95"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.927.53"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "5.530"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "6.669"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "91.67"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "0.000008869"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "28.285.52"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").Value = "14.93"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "5.175"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "2.182.49"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("D24").Value = "10.98"
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("D25").Value = "156.52"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "1.929"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "2.108"
$ws.Range("E28").Value = "  +4.46%  "
$ws.Range("D29").Value = "116.75"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "4.995"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "0.08935"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "3.359"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "1.255"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("D34").Value = "0.7758"
$ws.Range("E34").Value = "  +3.63%  "
$ws.Range("D35").Value = "4.716"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").Value = "2.707"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "1.110"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "0.5582"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").Value = "0.05339"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "3.011"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "7.071"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "0.1535"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "8.514"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "10.74"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "0.4894"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "107.15"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "1.669"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "69.09"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").Value = "0.06139"
$ws.Range("E51").Value = "  +0.53%  "
